# "new rates as of 3/18"
# Updates the rate table: new rates for Morris, IL and Belleville, MI,
# a new weight-break column (15 / column P) for every lane, and three
# brand-new destination rows (Monroe Township NJ, Monrovia MD, Owatonna MN).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number format used by every rate cell in the table (style index 1 in the
# original workbook: numFmtId 2, i.e. "0.00").
$rateFormat = "0.00"

function Set-Rate($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = $rateFormat
    if ($null -ne $value) {
        $cell.Value = $value
    }
}

# ---------------------------------------------------------------------
# Header row: add the 15th weight-break column (P1 = 15)
# ---------------------------------------------------------------------
$ws.Range("P1").Value = 15

# ---------------------------------------------------------------------
# Row 2 - Brookshire, TX: values unchanged, just extend with empty P2
# ---------------------------------------------------------------------
Set-Rate 2 16 $null

# ---------------------------------------------------------------------
# Row 3 - Morris, IL: brand new rates, now populated through column P
# ---------------------------------------------------------------------
$row3 = @(290,570,850,1120,1400,1680,1945,2200,2475,2750,3000,3240,3510,3780,4015)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = 2 + $i
    Set-Rate 3 $col $row3[$i]
}

# ---------------------------------------------------------------------
# Row 4 - Sumner, WA: values unchanged, just extend with empty P4
# ---------------------------------------------------------------------
Set-Rate 4 16 $null

# ---------------------------------------------------------------------
# Row 5 - Tracy, CA: values unchanged, just extend with empty P5
# ---------------------------------------------------------------------
Set-Rate 5 16 $null

# ---------------------------------------------------------------------
# Row 13 - new destination: Monroe Township, NJ
# (string added ahead of "Belleville,  MI" so the shared-string table
# comes out in the same order the workbook author produced it in)
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "Monroe Township, NJ"
$row13 = @(470,880,1245,1580,1925,2280,2625,2920,3240,3550)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $col = 2 + $i
    Set-Rate 13 $col $row13[$i]
}
for ($c = 12; $c -le 16; $c++) {
    Set-Rate 13 $c $null
}

# ---------------------------------------------------------------------
# Row 14 - new destination: Monrovia, MD
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Monrovia, MD"
$row14 = @(470,880,1245,1580,1925,2280,2625,2920,3240,3550)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $col = 2 + $i
    Set-Rate 14 $col $row14[$i]
}
for ($c = 12; $c -le 16; $c++) {
    Set-Rate 14 $c $null
}

# ---------------------------------------------------------------------
# Row 6 - Belleville, MI (renamed "Belleville,  MI" with a double space):
# brand new rates, now populated through column P
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Belleville,  MI"
$row6 = @(465,800,1185,1540,1825,2160,2415,2760,3060,3400,3740,4080,4355,4690,5025)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $col = 2 + $i
    Set-Rate 6 $col $row6[$i]
}

# ---------------------------------------------------------------------
# Rows 7-11: values unchanged, just extend each with an empty column P
# ---------------------------------------------------------------------
Set-Rate 7 16 $null
Set-Rate 8 16 $null
Set-Rate 9 16 $null
Set-Rate 10 16 $null
Set-Rate 11 16 $null

# ---------------------------------------------------------------------
# Row 12 - Grand Prairie, TX: values unchanged (B-D), now extended with
# empty, formatted cells from E through P
# ---------------------------------------------------------------------
for ($c = 5; $c -le 16; $c++) {
    Set-Rate 12 $c $null
}

# ---------------------------------------------------------------------
# Row 15 - new destination: Owatonna, MN
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Owatonna, MN"
$row15 = @(525,1135,1475,1580,1840,2190,2345,2600,2880,3130)
for ($i = 0; $i -lt $row15.Length; $i++) {
    $col = 2 + $i
    Set-Rate 15 $col $row15[$i]
}
for ($c = 12; $c -le 16; $c++) {
    Set-Rate 15 $c $null
}

# ---------------------------------------------------------------------
# Cosmetics: column A widens to fit the longer destination names, and the
# final selection lands on the newly-added rows 15:16.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 20.5
$ws.Range("A15:XFD16").Select() | Out-Null
